$d = $word.ActiveDocument

# 1) Remove the whole "Meta description: ..." paragraph that currently sits
#    right after the "Play Diamond Digger for Free - Slot Game Review"
#    heading at the top of the document.
$findRange = $d.Content
$found = $findRange.Find.Execute("Meta description", $true, $false, $false, $false, $false, `
                                  $true, 1, $false, "", 0)
if ($found) {
    $metaPara = $findRange.Paragraphs(1)
    $metaPara.Range.Delete()
}

# 2) At the end of the document, the italic image-generation-prompt paragraph
#    ("Create a feature image...") is replaced by two new paragraphs: a bold
#    title line, followed by the text that used to be the meta description,
#    now rendered in italics.
$count = $d.Paragraphs.Count
$lastPara = $d.Paragraphs($count)

$xml = "<w:p xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'>" + `
       "<w:r/><w:r><w:rPr><w:b/></w:rPr>" + `
       "<w:t>Play Diamond Digger for Free - Slot Game Review</w:t></w:r></w:p>" + `
       "<w:p xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'>" + `
       "<w:r/><w:r><w:rPr><w:i/></w:rPr>" + `
       "<w:t>Read our review of Diamond Digger, a fun and easy-to-use slot game. Play for free and learn about the game's features and drawbacks.</w:t></w:r></w:p>"

$lastPara.Range.InsertXML($xml)
